$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A7:FI8").Copy($ws.Range("A9:FI10"))
